$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($addr in @('D5', 'D7', 'D10', 'D13', 'D14', 'D16', 'D19', 'D20', 'D22', 'D25', 'D26', 'D27', 'D29', 'D32', 'D34', 'D35', 'D36', 'D37', 'D38', 'D40', 'D42', 'D43', 'D44', 'D45', 'D46', 'D49', 'D50')) {
    $ws.Range($addr).NumberFormat = '@'
}

$ws.Range('D2').Value = '37.923.96'
$ws.Range('E2').Value = '  +0.38%  '

$ws.Range('D3').Value = '2.087.26'
$ws.Range('E3').Value = '  +2.99%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').Value = '228.35'
$ws.Range('E5').Value = '  +0.75%  '

$ws.Range('E6').Value = '  -0.45%  '

$ws.Range('D7').Value = '60.57'
$ws.Range('E7').Value = '  +1.43%  '

$ws.Range('E8').Value = '  +0.00%  '

$ws.Range('E9').Value = '  -0.08%  '

$ws.Range('D10').Value = '0.0839'
$ws.Range('E10').Value = '  +2.73%  '

$ws.Range('E11').Value = '  -0.01%  '

$ws.Range('D12').Value = '2.395.07'
$ws.Range('E12').Value = '  +2.89%  '

$ws.Range('D13').Value = '14.54'
$ws.Range('E13').Value = '  +0.89%  '

$ws.Range('D14').Value = '22.04'
$ws.Range('E14').Value = '  +5.39%  '

$ws.Range('E15').Value = '  +6.71%  '

$ws.Range('D16').Value = '0.769'
$ws.Range('E16').Value = '  +1.98%  '

$ws.Range('D17').Value = '2.089.69'
$ws.Range('E17').Value = '  +3.13%  '

$ws.Range('D18').Value = '37.903.32'
$ws.Range('E18').Value = '  +0.54%  '

$ws.Range('D19').Value = '6.05'
$ws.Range('E19').Value = '  +2.13%  '

$ws.Range('D20').Value = '69.88'
$ws.Range('E20').Value = '  +0.27%  '

$ws.Range('E21').Value = '  +1.55%  '

$ws.Range('D22').Value = '223.07'
$ws.Range('E22').Value = '  -0.46%  '

$ws.Range('E23').Value = '  +0.59%  '

$ws.Range('E24').Value = '  -0.29%  '

$ws.Range('D25').Value = '2.32'
$ws.Range('E25').Value = '  +3.38%  '

$ws.Range('D26').Value = '169.45'
$ws.Range('E26').Value = '  +2.22%  '

$ws.Range('D27').Value = '9.41'
$ws.Range('E27').Value = '  +1.92%  '

$ws.Range('E28').Value = '  +2.17%  '

$ws.Range('D29').Value = '18.96'
$ws.Range('E29').Value = '  -0.33%  '

$ws.Range('E30').Value = '  +3.41%  '

$ws.Range('E31').Value = '  -0.74%  '

$ws.Range('D32').Value = '2.35'
$ws.Range('E32').Value = '  +7.39%  '

$ws.Range('E33').Value = '  +0.89%  '

$ws.Range('D34').Value = '4.64'
$ws.Range('E34').Value = '  +3.95%  '

$ws.Range('D35').Value = '0.0606'
$ws.Range('E35').Value = '  +0.66%  '

$ws.Range('B36').Value = 'THORChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D36').Value = '6.48'
$ws.Range('E36').Value = '  +1.84%  '

$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').Value = '2.40'
$ws.Range('E37').Value = '  +6.81%  '

$ws.Range('D38').Value = '3.53'
$ws.Range('E38').Value = '  +9.64%  '

$ws.Range('E39').Value = '  +0.02%  '

$ws.Range('D40').Value = '18.22'
$ws.Range('E40').Value = '  +8.72%  '

$ws.Range('D41').Value = '1.544.30'
$ws.Range('E41').Value = '  +1.64%  '

$ws.Range('D42').Value = '100.04'
$ws.Range('E42').Value = '  +4.91%  '

$ws.Range('D43').Value = '0.0218'
$ws.Range('E43').Value = '  +0.81%  '

$ws.Range('D44').Value = '2.83'
$ws.Range('E44').Value = '  -0.33%  '

$ws.Range('D45').Value = '0.0906'
$ws.Range('E45').Value = '  -0.68%  '

$ws.Range('D46').Value = '4.14'
$ws.Range('E46').Value = '  +3.27%  '

$ws.Range('E47').Value = '  +1.63%  '

$ws.Range('E48').Value = '  +2.11%  '

$ws.Range('D49').Value = '2.99'
$ws.Range('E49').Value = '  +0.96%  '

$ws.Range('D50').Value = '7.18'
$ws.Range('E50').Value = '  +1.58%  '

$ws.Range('D51').Value = '2.282.22'
$ws.Range('E51').Value = '  +2.95%  '
